$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clean up trailing/doubled whitespace in the "Formula" column (col G) ---
# Row 2 (Management Fees): trailing tab+newline removed
$ws.Range("G2").Value = "capital_commitment.committed_amount_cents * ( fund_unit_setting.management_fee / 4 ) / 100"

# Row 4 (Unit Premium): double space before "/ 100" collapsed, trailing tab removed
$ws.Range("G4").Value = "capital_commitment.fund.total_units_premium_cents * capital_commitment.percentage / 100"

# Row 5 (Generate Investable Capital): double spaces collapsed
$ws.Range("G5").Value = 'capital_commitment["properties"]["investable_capital"] = ( capital_commitment.collected_amount_cents - capital_commitment.cumulative_account_entry("Management Fees").amount_cents + capital_commitment.cumulative_account_entry("Unit Premium").amount_cents ) / 100'

# Row 7 (Accured Interest): double space collapsed, trailing tab+newline removed
$ws.Range("G7").Value = "fund_account_entry.amount_cents * capital_commitment.properties['investable_capital_percentage'] / 100.0"

# --- Update the saved view state: active cell / scroll position moved from I10 to G5 ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 7
$win.ScrollRow = 1
$ws.Range("G5").Select() | Out-Null
